# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3030
    4  = 218
    6  = 193
    7  = 1649
    8  = 1615
    13 = 187
    16 = 232
    20 = 39
    26 = 2031
    28 = 459
    30 = 179
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
